$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,(2000, 585400, 556400)
    ,(823, 5500, 7300)
    ,(750, 2800, 5400)
    ,(1571, 2825, 5750)
    ,(412, 11400, 21700)
    ,(385, 3000, 5100)
    ,(1027, 11600, 13700)
    ,(750, 2675, 4366)
    ,(792, 11500, 10050)
    ,(385, 2700, 5200)
    ,(392, 3466, 5700)
    ,(318, 2725, 3733)
    ,(815, 11200, 9100)
    ,(1900, 2725, 10800)
    ,(680, 5250, 6150)
    ,(623, 5400, 5850)
    ,(542, 7233, 6000)
    ,(1300, 4600, 6100)
    ,(1168, 5350, 6100)
    ,(4774950, 11700, 6550)
    ,(1636, 5400, 6150)
    ,(385, 8466, 6300)
    ,(412, 7500, 10200)
    ,(364, 3400, 7650)
    ,(452, 706, 13100)
    ,(346, 742, 14850)
    ,(945, 784, 6250)
    ,(412, 675, 6900)
    ,(3044, 785, 10400)
    ,(355, 680, 6700)
    ,(375, 981, 6000)
    ,(400, 866, 13400)
    ,(358, 568, 8100)
    ,(9650, 617, 6400)
    ,(778, 600, 26200)
    ,(392, 800, 6400)
    ,(2000, 588, 10200)
    ,(278, 792, 5500)
    ,(1097, 547, 7650)
    ,(166, 662, 11200)
    ,(468, 611, 10600)
    ,(396, 656, 11950)
    ,(351, 600, 5850)
    ,(206, 784, 10600)
    ,(177, 771, 6400)
    ,(198, 2200, 6450)
    ,(194, 1163, 17800)
    ,(231, 958, 8500)
    ,(245, 605, 6000)
    ,(204, 706, 6400)
    ,(214, 568, 5650)
    ,(187, 750, 8200)
    ,(226, 1275, 31050)
    ,(231, 757, 5100)
    ,(210, 1155, 18800)
    ,(230, 680, 1783)
    ,(192, 3400, 2100)
    ,(214, 311, 4033)
    ,(196, 300, 1866)
    ,(223, 364, 2625)
    ,(196, 367, 2280)
    ,(192, 328, 2100)
    ,(183, 321, 2100)
    ,(174, 500, 2180)
    ,(256, 351, 2800)
    ,(275, 355, 3150)
    ,(174, 288, 1966)
    ,(196, 477, 3633)
    ,(242, 308, 3120)
    ,(210, 335, 2260)
    ,(212, 221, 1155)
    ,(114, 185, 1030)
    ,(152, 623, 2750)
    ,(245, 332, 2200)
    ,(210, 408, 2200)
    ,(294, 510, 2380)
    ,(188, 271, 2080)
    ,(210, 385, 2040)
    ,(170, 245, 2700)
    ,(118, 176, 1733)
    ,(110, 515, 6850)
    ,(117, 219, 1275)
    ,(114, 174, 735)
    ,(112, 208, 500)
    ,(113, 343, 500)
    ,(204, 358, 1471)
    ,(214, 425, 936)
    ,(201, 329, 900)
    ,(182, 350, 990)
    ,(118, 214, 510)
    ,(143, 178, 566)
    ,(217, 1700, 936)
    ,(210, 201, 963)
    ,(206, 245, 972)
    ,(429, 237, 1009)
    ,(196, 212, 945)
    ,(212, 247, 2540)
    ,(196, 247, 954)
    ,(201, 245, 981)
    ,(258, 236, 866)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowVals[0]
    $ws.Cells.Item($row, 2).Value = $rowVals[1]
    $ws.Cells.Item($row, 3).Value = $rowVals[2]
}